# Append three new blank ("Blank" layout) slides to the end of the
# sketchbook deck, mirroring the plain blank slides already used
# throughout the deck (e.g. slide 20): empty shape tree, no legacy
# vmlDrawing relationship, just a link to the Blank slide layout.

$p = $ppt.ActivePresentation

# Use the layout already applied to the last existing slide (the
# "Blank" layout) as the template layout for the new slides.
$blankLayout = $p.Slides.Item($p.Slides.Count).CustomLayout

for ($n = 0; $n -lt 3; $n++) {
    # Duplicating the title slide (slide 1) gives us a fresh slide that
    # already carries the normal PowerPoint boilerplate (xfrm/clrMapOvr/
    # timing) without dragging along the stray vmlDrawing relationship
    # that the ink/drawing slides in this deck have.
    $newSlide = $p.Slides.Item(1).Duplicate().Item(1)

    # Clear out the inherited placeholder/text shapes so the slide ends
    # up with a totally empty shape tree. Placeholders need an extra
    # delete pass since the first pass just resets them to the layout
    # default instead of removing them outright.
    for ($i = $newSlide.Shapes.Count; $i -ge 1; $i--) {
        $newSlide.Shapes.Item($i).Delete()
    }
    for ($i = $newSlide.Shapes.Count; $i -ge 1; $i--) {
        $newSlide.Shapes.Item($i).Delete()
    }

    # Swap the layout from "Title Slide" to "Blank".
    $newSlide.CustomLayout = $blankLayout

    # Move the freshly duplicated slide (currently right after slide 1)
    # to the end of the deck.
    $newSlide.MoveTo($p.Slides.Count)
}
